# Final version. Ready to merge
#
# Rescale the "rate" columns (F/G, J/K, O/P, T/U, Y/Z) on Sheet1 from their
# small calculated/literal values to round placeholder numbers, and update
# the sheet's scroll position / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 block: F/G columns (rows 2-10 only) ---
$ws.Range("F2:F10").Value = 20
$ws.Range("G2:G10").Value = 20

# J/K were formulas (=5.2*0.8 / =5.2*0.6, shared across J3:J10 / K3:K10) -
# overwrite with plain literals, which clears the formulas.
$ws.Range("J2:J10").Value = 30
$ws.Range("K2:K10").Value = 30

# --- O/P columns (rows 2-39) ---
$ws.Range("O2:O39").Value = 21
$ws.Range("P2:P39").Value = 21

# T/U were formulas (=5.2*0.8 / =5.2*0.6, shared across T3:T39 / U3:U39)
$ws.Range("T2:T39").Value = 22
$ws.Range("U2:U39").Value = 22

# Y/Z were formulas (=8.1*0.8 / =8.1*0.6, shared across Y3:Y39 / Z3:Z39)
$ws.Range("Y2:Y39").Value = 23
$ws.Range("Z2:Z39").Value = 23

# --- View state: scroll so column H is the left-most visible column, and
# select Y2:Z39 (active cell Y2) ---
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Y2:Z39").Select()
